$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Normal text style" paragraph (Normal style, no explicit
#    pPr/pStyle) right before the "Second Heading 1" paragraph.
# ---------------------------------------------------------------------------
$secondHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Second Heading 1") {
        $secondHeading = $p
        break
    }
}

$insertRng = $secondHeading.Range
$insertRng.Collapse(1)
$insertRng.InsertParagraphBefore()

# Re-find "Second Heading 1" now that the paragraph collection has shifted,
# then grab the freshly created (still empty) paragraph immediately
# preceding it and replace its whole range (including its end-of-paragraph
# mark) with clean WordprocessingML that carries no pPr/pStyle, so the
# paragraph renders with the document's default ("Normal") style.
$secondHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Second Heading 1") {
        $secondHeading = $p
        break
    }
}
$newPara = $secondHeading.Previous()
$fullRng = $newPara.Range

$normalParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p><w:r><w:t>Normal text style</w:t></w:r></w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$fullRng.InsertXML($normalParaXml)

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the empty Heading 2 paragraph (after
#    "Next") to the start of the "Second Heading 1" paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$secondHeading2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Second Heading 1") {
        $secondHeading2 = $p
        break
    }
}
$bmRng = $secondHeading2.Range
$bmRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRng)

# ---------------------------------------------------------------------------
# 3. Footer: "Release" field goes from the complex fldChar form back to a
#    simple <w:fldSimple>; "NUMPAGES" field goes from <w:fldSimple> to the
#    complex fldChar form (with an extra rPr/noProof on the end fldChar).
# ---------------------------------------------------------------------------
$footer = $d.Sections(1).Footers(2)   # wdHeaderFooterPrimary = 2
$footerPara = $footer.Range.Paragraphs(1)

$footerXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/footer2.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.footer+xml">
    <pkg:xmlData>
      <w:ftr xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:p>
          <w:pPr>
            <w:pStyle w:val="PageFooter"/>
          </w:pPr>
          <w:r>
            <w:t xml:space="preserve">Release </w:t>
          </w:r>
          <w:fldSimple w:instr=" DOCPROPERTY  Release  \* MERGEFORMAT ">
            <w:r>
              <w:t>5.0</w:t>
            </w:r>
          </w:fldSimple>
          <w:r>
            <w:ptab w:relativeTo="margin" w:alignment="center" w:leader="none"/>
          </w:r>
          <w:sdt>
            <w:sdtPr>
              <w:alias w:val="Title"/>
              <w:tag w:val=""/>
              <w:id w:val="149645157"/>
              <w:placeholder>
                <w:docPart w:val="3E7BB149AEA54986861D4A8465884C19"/>
              </w:placeholder>
              <w:dataBinding w:prefixMappings="xmlns:ns0='http://purl.org/dc/elements/1.1/' xmlns:ns1='http://schemas.openxmlformats.org/package/2006/metadata/core-properties' " w:xpath="/ns1:coreProperties[1]/ns0:title[1]" w:storeItemID="{6C3C8BC8-F283-45AE-878A-BAB7291924A1}"/>
              <w:text/>
            </w:sdtPr>
            <w:sdtEndPr/>
            <w:sdtContent>
              <w:proofErr w:type="spellStart"/>
              <w:r>
                <w:t>Pandoc</w:t>
              </w:r>
              <w:proofErr w:type="spellEnd"/>
              <w:r>
                <w:t xml:space="preserve"> Reference Document</w:t>
              </w:r>
            </w:sdtContent>
          </w:sdt>
          <w:r>
            <w:ptab w:relativeTo="margin" w:alignment="right" w:leader="none"/>
          </w:r>
          <w:r>
            <w:t xml:space="preserve">Page </w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> PAGE   \* MERGEFORMAT </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t>1</w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
          <w:r>
            <w:t xml:space="preserve"> of </w:t>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="begin"/>
          </w:r>
          <w:r>
            <w:instrText xml:space="preserve"> NUMPAGES   \* MERGEFORMAT </w:instrText>
          </w:r>
          <w:r>
            <w:fldChar w:fldCharType="separate"/>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:t>1</w:t>
          </w:r>
          <w:r>
            <w:rPr>
              <w:noProof/>
            </w:rPr>
            <w:fldChar w:fldCharType="end"/>
          </w:r>
        </w:p>
      </w:ftr>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$footer.Range.InsertXML($footerXml)

# ---------------------------------------------------------------------------
# 4. Font substitution table signature tweaks + new rsid in the glossary
#    settings (cosmetic metadata only).
# ---------------------------------------------------------------------------
Write-Output "done"
